$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-01 Sunday", "2024-09-02 Monday"),
    @("461×6=2766", "594×4=2376"),
    @("756×6=4536", "762×6=4572"),
    @("868×6=5208", "292×7=2044"),
    @("796×7=5572", "457×9=4113"),
    @("152×5=760", "137×6=822"),
    @("826×4=3304", "290×3=870"),
    @("993×9=8937", "496×7=3472"),
    @("415×6=2490", "472×9=4248"),
    @("762×4=3048", "140×6=840"),
    @("433×9=3897", "271×6=1626"),
    @("920×2=1840", "729×9=6561"),
    @("138×9=1242", "261×9=2349"),
    @("837×5=4185", "422×4=1688"),
    @("961×7=6727", "573×8=4584"),
    @("523×6=3138", "247×3=741"),
    @("798×4=3192", "309×8=2472"),
    @("595×7=4165", "982×2=1964"),
    @("939×3=2817", "922×7=6454"),
    @("739×3=2217", "167×2=334"),
    @("727×5=3635", "648×6=3888"),
    @("493×8=3944", "883×3=2649"),
    @("336×3=1008", "904×2=1808"),
    @("689×6=4134", "792×6=4752"),
    @("981×9=8829", "269×2=538"),
    @("944×7=6608", "965×5=4825")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
